# Scheduled-runner refresh of FFXIV market-board data (currentAveragePrice*,
# LevePrice*, LeveProfit*) across the per-job worksheets. Plain values only —
# no formulas in these columns.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# row 2
$ws.Range("H2").Value = 1419.8
$ws.Range("I2").Value = 300
$ws.Range("J2").Value = 2166.3333
$ws.Range("K2").Value = 300
$ws.Range("L2").Value = 2166.3333
$ws.Range("M2").Value = -187
$ws.Range("N2").Value = -2392.3333

# row 9
$ws.Range("H9").Value = 4871.5415
$ws.Range("I9").Value = 7188.067
$ws.Range("J9").Value = 1010.6667
$ws.Range("K9").Value = 7188.067
$ws.Range("L9").Value = 1010.6667
$ws.Range("M9").Value = -7019.067
$ws.Range("N9").Value = -1348.6667

# row 40
$ws.Range("H40").Value = 6955255
$ws.Range("J40").Value = 22247020
$ws.Range("L40").Value = 22247020
$ws.Range("N40").Value = -22247370

# row 43
$ws.Range("H43").Value = 4332.6665
$ws.Range("I43").Value = 3999.5
$ws.Range("J43").Value = 4999
$ws.Range("K43").Value = 3999.5
$ws.Range("L43").Value = 4999
$ws.Range("M43").Value = -3930.5
$ws.Range("N43").Value = -5137

# row 86
$ws.Range("H86").Value = 4314.6665
$ws.Range("I86").Value = 4000
$ws.Range("J86").Value = 4472
$ws.Range("K86").Value = 4000
$ws.Range("L86").Value = 4472
$ws.Range("M86").Value = -2877
$ws.Range("N86").Value = -6718

# row 89
$ws.Range("H89").Value = 4314.6665
$ws.Range("I89").Value = 4000
$ws.Range("J89").Value = 4472
$ws.Range("K89").Value = 20000
$ws.Range("L89").Value = 22360
$ws.Range("M89").Value = -14384
$ws.Range("N89").Value = -33592

# row 103
$ws.Range("H103").Value = 2010.8
$ws.Range("I103").Value = 1777.5
$ws.Range("J103").Value = 2166.3333
$ws.Range("K103").Value = 5332.5
$ws.Range("L103").Value = 6498.999899999999
$ws.Range("M103").Value = -4746.5
$ws.Range("N103").Value = -7670.999899999999

# row 106
$ws.Range("H106").Value = 3427.0625
$ws.Range("I106").Value = 2778
$ws.Range("K106").Value = 2778
$ws.Range("M106").Value = -2147

# row 132
$ws.Range("H132").Value = 2978
$ws.Range("I132").Value = 3063.1738
$ws.Range("J132").Value = 1998.5
$ws.Range("K132").Value = 9189.5214
$ws.Range("L132").Value = 5995.5
$ws.Range("M132").Value = -6659.5214
$ws.Range("N132").Value = -11055.5

# row 137
$ws.Range("H137").Value = 2021.2667
$ws.Range("I137").Value = 1316.3636
$ws.Range("J137").Value = 3959.75
$ws.Range("K137").Value = 3949.0908
$ws.Range("L137").Value = 11879.25
$ws.Range("M137").Value = -1399.0908
$ws.Range("N137").Value = -16979.25

# row 138
$ws.Range("H138").Value = 2782.6038
$ws.Range("J138").Value = 3218.3333
$ws.Range("L138").Value = 9654.999899999999
$ws.Range("N138").Value = -19934.9999

$ws = $wb.Worksheets.Item("ARM")
# row 2
$ws.Range("H2").Value = 3367.625
$ws.Range("I2").Value = 2562
$ws.Range("J2").Value = 5784.5
$ws.Range("K2").Value = 2562
$ws.Range("L2").Value = 5784.5
$ws.Range("M2").Value = -2449
$ws.Range("N2").Value = -6010.5

# row 32
$ws.Range("H32").Value = 4274.4604
$ws.Range("I32").Value = 2088.1633
$ws.Range("K32").Value = 2088.1633
$ws.Range("M32").Value = -1801.1633

# row 61
$ws.Range("H61").Value = 21742532
$ws.Range("I61").Value = 22730648
$ws.Range("K61").Value = 22730648
$ws.Range("M61").Value = -22730436

# row 74
$ws.Range("H74").Value = 41671196
$ws.Range("I74").Value = 76929170
$ws.Range("J74").Value = 2681.818
$ws.Range("K74").Value = 76929170
$ws.Range("L74").Value = 2681.818
$ws.Range("M74").Value = -76928296
$ws.Range("N74").Value = -4429.818

# row 77
$ws.Range("H77").Value = 41671196
$ws.Range("I77").Value = 76929170
$ws.Range("J77").Value = 2681.818
$ws.Range("K77").Value = 384645850
$ws.Range("L77").Value = 13409.09
$ws.Range("M77").Value = -384641482
$ws.Range("N77").Value = -22145.09

# row 86
$ws.Range("H86").Value = 49499.5
$ws.Range("J86").Value = 49499.5
$ws.Range("L86").Value = 49499.5
$ws.Range("N86").Value = -51871.5

# row 89
$ws.Range("H89").Value = 49499.5
$ws.Range("J89").Value = 49499.5
$ws.Range("L89").Value = 148498.5
$ws.Range("N89").Value = -160354.5

# row 116
$ws.Range("H116").Value = 3367.625
$ws.Range("I116").Value = 2562
$ws.Range("J116").Value = 5784.5
$ws.Range("K116").Value = 2562
$ws.Range("L116").Value = 5784.5
$ws.Range("M116").Value = -268
$ws.Range("N116").Value = -10372.5

# row 117
$ws.Range("H117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("L117").Value = 0
$ws.Range("N117").ClearContents()

# row 122
$ws.Range("H122").Value = 1671.2354
$ws.Range("I122").Value = 1750.75
$ws.Range("J122").Value = 399
$ws.Range("K122").Value = 5252.25
$ws.Range("L122").Value = 1197
$ws.Range("M122").Value = -2802.25
$ws.Range("N122").Value = -6097

# row 132
$ws.Range("H132").Value = 3035373.5
$ws.Range("I132").Value = 3230834.5
$ws.Range("J132").Value = 5724.5
$ws.Range("K132").Value = 9692503.5
$ws.Range("L132").Value = 17173.5
$ws.Range("M132").Value = -9689973.5
$ws.Range("N132").Value = -22233.5

# row 136
$ws.Range("H136").Value = 21742532
$ws.Range("I136").Value = 22730648
$ws.Range("K136").Value = 68191944
$ws.Range("M136").Value = -68189394

$ws = $wb.Worksheets.Item("BSM")
# row 3
$ws.Range("H3").Value = 3367.625
$ws.Range("I3").Value = 2562
$ws.Range("J3").Value = 5784.5
$ws.Range("K3").Value = 2562
$ws.Range("L3").Value = 5784.5
$ws.Range("M3").Value = -2448
$ws.Range("N3").Value = -6012.5

# row 134
$ws.Range("H134").Value = 11113665
$ws.Range("I134").Value = 12822354
$ws.Range("K134").Value = 38467062
$ws.Range("M134").Value = -38464527

$ws = $wb.Worksheets.Item("CRP")
# row 86
$ws.Range("H86").Value = 7699
$ws.Range("I86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("M86").ClearContents()

# row 88
$ws.Range("H88").Value = 15758.714
$ws.Range("J88").Value = 15758.714
$ws.Range("L88").Value = 15758.714
$ws.Range("N88").Value = -16570.714

# row 89
$ws.Range("H89").Value = 7699
$ws.Range("I89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("M89").ClearContents()

# row 91
$ws.Range("H91").Value = 15758.714
$ws.Range("J91").Value = 15758.714
$ws.Range("L91").Value = 15758.714
$ws.Range("N91").Value = -18566.714

# row 122
$ws.Range("H122").Value = 1642.0394
$ws.Range("I122").Value = 1727.0328
$ws.Range("J122").Value = 1296.4
$ws.Range("K122").Value = 5181.0984
$ws.Range("L122").Value = 3889.2
$ws.Range("M122").Value = -2731.0984
$ws.Range("N122").Value = -8789.200000000001

# row 132
$ws.Range("H132").Value = 100000970
$ws.Range("I132").Value = 125000960
$ws.Range("J132").Value = 999
$ws.Range("K132").Value = 375002880
$ws.Range("L132").Value = 2997
$ws.Range("M132").Value = -375000350
$ws.Range("N132").Value = -8057

$ws = $wb.Worksheets.Item("CUL")
# row 88
$ws.Range("H88").Value = 11047.429
$ws.Range("I88").Value = 4500
$ws.Range("J88").Value = 12138.667
$ws.Range("K88").Value = 13500
$ws.Range("L88").Value = 36416.001
$ws.Range("M88").Value = -13072
$ws.Range("N88").Value = -37272.001

# row 91
$ws.Range("H91").Value = 11047.429
$ws.Range("I91").Value = 4500
$ws.Range("J91").Value = 12138.667
$ws.Range("K91").Value = 13500
$ws.Range("L91").Value = 36416.001
$ws.Range("M91").Value = -12018
$ws.Range("N91").Value = -39380.001

# row 132
$ws.Range("H132").Value = 1841.9524
$ws.Range("J132").Value = 1769.8572
$ws.Range("L132").Value = 15928.7148
$ws.Range("N132").Value = -20988.7148

$ws = $wb.Worksheets.Item("GSM")
# row 126
$ws.Range("H126").Value = 4665.8965
$ws.Range("I126").Value = 4665.8965
$ws.Range("K126").Value = 13997.6895
$ws.Range("M126").Value = -11527.6895

# row 132
$ws.Range("H132").Value = 8930248
$ws.Range("I132").Value = 10418290
$ws.Range("J132").Value = 2000
$ws.Range("K132").Value = 31254870
$ws.Range("L132").Value = 6000
$ws.Range("M132").Value = -31252340
$ws.Range("N132").Value = -11060

# row 140
$ws.Range("H140").Value = 86665
$ws.Range("J140").Value = 86665
$ws.Range("L140").Value = 86665
$ws.Range("N140").Value = -97025

# row 141
$ws.Range("H141").Value = 84614.664
$ws.Range("J141").Value = 84614.664
$ws.Range("L141").Value = 84614.664
$ws.Range("N141").Value = -94974.664

$ws = $wb.Worksheets.Item("LTW")
# row 50
$ws.Range("H50").Value = 30000
$ws.Range("J50").Value = 30000
$ws.Range("L50").Value = 30000
$ws.Range("N50").Value = -31274

# row 122
$ws.Range("H122").Value = 3405.6956
$ws.Range("I122").Value = 3339.5454
$ws.Range("J122").Value = 3466.3333
$ws.Range("K122").Value = 10018.6362
$ws.Range("L122").Value = 10398.9999
$ws.Range("M122").Value = -7568.636200000001
$ws.Range("N122").Value = -15298.9999

# row 132
$ws.Range("H132").Value = 40010090
$ws.Range("I132").Value = 53346256
$ws.Range("J132").Value = 1596
$ws.Range("K132").Value = 160038768
$ws.Range("L132").Value = 4788
$ws.Range("M132").Value = -160036238
$ws.Range("N132").Value = -9848

# row 136
$ws.Range("H136").Value = 2540.5
$ws.Range("I136").Value = 1831.8334
$ws.Range("J136").Value = 3249.1667
$ws.Range("K136").Value = 5495.5002
$ws.Range("L136").Value = 9747.500100000001
$ws.Range("M136").Value = -2945.5002
$ws.Range("N136").Value = -14847.5001

$ws = $wb.Worksheets.Item("WVR")
# row 54
$ws.Range("H54").Value = 0
$ws.Range("J54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("N54").ClearContents()
